$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.880.36'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.280.92'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.95%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.35'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.638'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.69'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +9.79%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.654'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +13.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.68'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0978'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.71'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.39'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +7.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.105'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.616.54'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.01'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.889'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.258.31'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.787.97'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.21%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.33'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.46'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.02'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.14'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.90'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.61'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.45'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.13'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.42'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.08'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.41'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +10.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.131'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +7.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0804'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '30.77'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +25.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.127'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.68'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +18.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.77'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.20%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.33'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.26'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +16.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.91'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.212'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +12.17%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.20'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +8.14%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.98'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '61.69'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.23%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.80%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.25%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.20'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.88%  '
